# Automatische test-sync: 2025-08-01 23:38:50
#
# Adds a second test-mail row to the "Logs" sheet, a matching tally row
# to the "Dashboard" sheet, extends the conditional-formatting ranges to
# cover the new row and widens the bar chart's category/value series to
# include the new Dashboard row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 7 with the new test-mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A7").Value = "Wil je dit oppakken?"
$logs.Range("B7").Value = "mailmind.test@zohomail.eu"
$logs.Range("C7").Value = "Testmail #2: Wil je dit oppakken?"
$logs.Range("D7").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E7").Value = "Beste,`nDank voor je bericht. Kun je wat meer context geven over wat je precies bedoelt met `"Testmail #2`" en wat er van mij verwacht wordt om op te pakken? Graag hoor ik meer details, zodat ik je beter van dienst kan zijn.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F7").Value = "2025-08-01 23:38:18"
$logs.Range("G7").Value = "Ja"
$logs.Range("H7").Value = "Nee"
$logs.Range("I7").Value = "Ja"
$logs.Range("J7").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Logs sheet: extend conditional formatting ranges from row 6 to 7
# ---------------------------------------------------------------------
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range($col + "2:" + $col + "6")
    $newRange = $logs.Range($col + "2:" + $col + "7")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. Dashboard sheet: append row 3 with the new category tally
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B3").Value = 1

# ---------------------------------------------------------------------
# 4. Dashboard chart: widen the series' category/value references
#    so they cover the newly added row (A2:A3 / B2:B3)
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$3,'Dashboard'!`$B`$2:`$B`$3,1)"
